# Add a helper "index" column (I) to the dictionary sheet, numbering
# data rows 1..203 (row 1 is the header row and is left untouched).
# This supports the new "search by substring, group and subgroup" feature
# described in the commit message, which needs a stable row index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$lastRow = 204

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)   # column I = 9
    $cell.Value = $r - 1
    $cell.Style = "Normal"          # match column I's own number style (s="6"),
                                     # not any row-level custom format (s="8")
}

# Reflect the selection/scroll position that Excel recorded after this
# edit: column I (the newly filled range) selected.
$ws.Range("I2:I204").Select() | Out-Null
